$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-11-19 Wednesday" "2025-11-20 Thursday"

Replace-Text "49×21=" "41×15="
Replace-Text "54×74=" "23×14="
Replace-Text "77×26=" "48×63="
Replace-Text "89×62=" "91×41="
Replace-Text "17×79=" "53×53="

Replace-Text "58×34=" "60×64="
Replace-Text "97×70=" "84×18="
Replace-Text "47×70=" "24×29="
Replace-Text "80×42=" "70×44="
Replace-Text "65×92=" "48×55="

Replace-Text "33×65=" "16×54="
Replace-Text "73×50=" "89×35="
Replace-Text "74×77=" "96×90="
Replace-Text "42×46=" "68×92="
Replace-Text "30×76=" "41×75="

Replace-Text "84×36=" "17×56="
Replace-Text "94×33=" "67×19="
Replace-Text "65×55=" "25×90="
Replace-Text "94×61=" "95×98="
Replace-Text "28×65=" "60×49="

Replace-Text "32×47=" "17×17="
Replace-Text "73×72=" "54×25="
Replace-Text "46×35=" "64×78="
Replace-Text "26×58=" "34×88="
Replace-Text "80×92=" "16×65="
